$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: swap "average_doctor" / "average_doctor_old" column titles ---
# (the "average_doctor_old" column now appears before "average_doctor")
$ws.Range("BP1").Value2 = "average_doctor_old"
$ws.Range("BQ1").Value2 = "average_doctor"

# --- Updated per-app "_old" statistics (recomputed after adding the Harvard case classification) ---

# Row 4
$ws.Range("E4").Value2 = 0.555
$ws.Range("F4").Value2 = 0.067
$ws.Range("G4").Value2 = 0.259
$ws.Range("N4").Value2 = 0.542
$ws.Range("O4").Value2 = 0.092
$ws.Range("P4").Value2 = 0.303
$ws.Range("Q4").Value2 = 0.157
$ws.Range("R4").Value2 = 0.089
$ws.Range("S4").Value2 = 0.299
$ws.Range("W4").Value2 = 0.298
$ws.Range("X4").Value2 = 0.108
$ws.Range("Y4").Value2 = 0.328
$ws.Range("AI4").Value2 = 0.192
$ws.Range("AJ4").Value2 = 0.079
$ws.Range("AK4").Value2 = 0.281
$ws.Range("AU4").Value2 = 0.202
$ws.Range("AV4").Value2 = 0.024
$ws.Range("AW4").Value2 = 0.154
$ws.Range("BA4").Value2 = 1.87
$ws.Range("BB4").Value2 = 0.132
$ws.Range("BC4").Value2 = 0.364
$ws.Range("BG4").Value2 = 0.655
$ws.Range("BH4").Value2 = 0.192
$ws.Range("BI4").Value2 = 0.439
$ws.Range("BM4").Value2 = 0.655
$ws.Range("BN4").Value2 = 0.109
$ws.Range("BO4").Value2 = 0.33
$ws.Range("BP4").Value2 = 0.623
$ws.Range("BQ4").Value2 = 0.702

# Row 5
$ws.Range("E5").Value2 = 0.5639999999999999
$ws.Range("F5").Value2 = 0.049
$ws.Range("G5").Value2 = 0.222
$ws.Range("N5").Value2 = 0.638
$ws.Range("O5").Value2 = 0.052
$ws.Range("P5").Value2 = 0.228
$ws.Range("Q5").Value2 = 0.104
$ws.Range("R5").Value2 = 0.036
$ws.Range("S5").Value2 = 0.19
$ws.Range("W5").Value2 = 0.211
$ws.Range("X5").Value2 = 0.076
$ws.Range("Y5").Value2 = 0.276
$ws.Range("AI5").Value2 = 0.201
$ws.Range("AJ5").Value2 = 0.081
$ws.Range("AK5").Value2 = 0.284
$ws.Range("AU5").Value2 = 0.364
$ws.Range("AV5").Value2 = 0.08799999999999999
$ws.Range("AW5").Value2 = 0.297
$ws.Range("BA5").Value2 = 1.073
$ws.Range("BB5").Value2 = 0.029
$ws.Range("BC5").Value2 = 0.17
$ws.Range("BG5").Value2 = 0.358
$ws.Range("BH5").Value2 = 0.062
$ws.Range("BI5").Value2 = 0.25
$ws.Range("BM5").Value2 = 0.365
$ws.Range("BN5").Value2 = 0.028
$ws.Range("BO5").Value2 = 0.166
$ws.Range("BP5").Value2 = 0.358
$ws.Range("BQ5").Value2 = 0.393

# Row 6
$ws.Range("E6").Value2 = 0.5590000000000001
$ws.Range("N6").Value2 = 0.586
$ws.Range("Q6").Value2 = 0.125
$ws.Range("W6").Value2 = 0.247
$ws.Range("AI6").Value2 = 0.196
$ws.Range("AU6").Value2 = 0.26
$ws.Range("BA6").Value2 = 1.363
$ws.Range("BG6").Value2 = 0.463
$ws.Range("BM6").Value2 = 0.469
$ws.Range("BP6").Value2 = 0.454
$ws.Range("BQ6").Value2 = 0.502

# Row 7
$ws.Range("E7").Value2 = 0.5620000000000001
$ws.Range("N7").Value2 = 0.616
$ws.Range("Q7").Value2 = 0.112
$ws.Range("W7").Value2 = 0.224
$ws.Range("AI7").Value2 = 0.199
$ws.Range("AU7").Value2 = 0.314
$ws.Range("BA7").Value2 = 1.172
$ws.Range("BG7").Value2 = 0.394
$ws.Range("BM7").Value2 = 0.4
$ws.Range("BP7").Value2 = 0.391
$ws.Range("BQ7").Value2 = 0.43

# Row 8
$ws.Range("E8").Value2 = 0.641
$ws.Range("F8").Value2 = 0.081
$ws.Range("G8").Value2 = 0.284
$ws.Range("N8").Value2 = 0.849
$ws.Range("O8").Value2 = 0.017
$ws.Range("P8").Value2 = 0.13
$ws.Range("Q8").Value2 = 0.114
$ws.Range("R8").Value2 = 0.065
$ws.Range("S8").Value2 = 0.256
$ws.Range("W8").Value2 = 0.341
$ws.Range("X8").Value2 = 0.126
$ws.Range("Y8").Value2 = 0.354
$ws.Range("AI8").Value2 = 0.208
$ws.Range("AJ8").Value2 = 0.111
$ws.Range("AK8").Value2 = 0.333
$ws.Range("AU8").Value2 = 0.327
$ws.Range("AV8").Value2 = 0.08799999999999999
$ws.Range("AW8").Value2 = 0.296
$ws.Range("BA8").Value2 = 1.633
$ws.Range("BB8").Value2 = 0.099
$ws.Range("BC8").Value2 = 0.314
$ws.Range("BG8").Value2 = 0.53
$ws.Range("BH8").Value2 = 0.138
$ws.Range("BI8").Value2 = 0.371
$ws.Range("BM8").Value2 = 0.61
$ws.Range("BN8").Value2 = 0.07099999999999999
$ws.Range("BO8").Value2 = 0.266
$ws.Range("BP8").Value2 = 0.544
$ws.Range("BQ8").Value2 = 0.573

# Row 9
$ws.Range("E9").Value2 = 0.571
$ws.Range("F9").Value2 = 0.245
$ws.Range("G9").Value2 = 0.495
$ws.Range("N9").Value2 = 0.929
$ws.Range("O9").Value2 = 0.066
$ws.Range("P9").Value2 = 0.258
$ws.Range("AI9").Value2 = 0.214
$ws.Range("AJ9").Value2 = 0.168
$ws.Range("AK9").Value2 = 0.41
$ws.Range("BA9").Value2 = 1.785
$ws.Range("BM9").Value2 = 0.714
$ws.Range("BN9").Value2 = 0.204
$ws.Range("BO9").Value2 = 0.452
$ws.Range("BP9").Value2 = 0.595
$ws.Range("BQ9").Value2 = 0.598

# Row 10
$ws.Range("E10").Value2 = 0.643
$ws.Range("F10").Value2 = 0.23
$ws.Range("G10").Value2 = 0.479
$ws.Range("N10").Value2 = 1
$ws.Range("O10").Value2 = 0
$ws.Range("P10").Value2 = 0
$ws.Range("W10").Value2 = 0.429
$ws.Range("X10").Value2 = 0.245
$ws.Range("Y10").Value2 = 0.495
$ws.Range("AI10").Value2 = 0.214
$ws.Range("AJ10").Value2 = 0.168
$ws.Range("AK10").Value2 = 0.41
$ws.Range("BA10").Value2 = 2
$ws.Range("BB10").Value2 = 0.245
$ws.Range("BC10").Value2 = 0.495
$ws.Range("BM10").Value2 = 0.786
$ws.Range("BN10").Value2 = 0.168
$ws.Range("BO10").Value2 = 0.41
$ws.Range("BP10").Value2 = 0.667
$ws.Range("BQ10").Value2 = 0.718

# Row 11
$ws.Range("E11").Value2 = 0.714
$ws.Range("F11").Value2 = 0.204
$ws.Range("G11").Value2 = 0.452
$ws.Range("N11").Value2 = 1
$ws.Range("O11").Value2 = 0
$ws.Range("P11").Value2 = 0
$ws.Range("W11").Value2 = 0.429
$ws.Range("X11").Value2 = 0.245
$ws.Range("Y11").Value2 = 0.495
$ws.Range("AI11").Value2 = 0.214
$ws.Range("AJ11").Value2 = 0.168
$ws.Range("AK11").Value2 = 0.41
$ws.Range("AU11").Value2 = 0.5
$ws.Range("AV11").Value2 = 0.25
$ws.Range("AW11").Value2 = 0.5
$ws.Range("BA11").Value2 = 2
$ws.Range("BB11").Value2 = 0.245
$ws.Range("BC11").Value2 = 0.495
$ws.Range("BM11").Value2 = 0.786
$ws.Range("BN11").Value2 = 0.168
$ws.Range("BO11").Value2 = 0.41
$ws.Range("BP11").Value2 = 0.667
$ws.Range("BQ11").Value2 = 0.718

# Row 12
$ws.Range("E12").Value2 = 1.6
$ws.Range("F12").Value2 = 1.64
$ws.Range("G12").Value2 = 1.281
$ws.Range("N12").Value2 = 1.071
$ws.Range("O12").Value2 = 0.066
$ws.Range("P12").Value2 = 0.258
$ws.Range("W12").Value2 = 1.167
$ws.Range("X12").Value2 = 0.139
$ws.Range("Y12").Value2 = 0.373
$ws.Range("AU12").Value2 = 3.5
$ws.Range("AV12").Value2 = 3.75
$ws.Range("AW12").Value2 = 1.936
$ws.Range("BA12").Value2 = 3.438
$ws.Range("BB12").Value2 = 0.109
$ws.Range("BC12").Value2 = 0.331
$ws.Range("BM12").Value2 = 1.091
$ws.Range("BN12").Value2 = 0.083
$ws.Range("BO12").Value2 = 0.287
$ws.Range("BP12").Value2 = 1.146
$ws.Range("BQ12").Value2 = 1.232

# Row 13
$ws.Range("BP13").Value2 = 0.6830000000000001
$ws.Range("BQ13").Value2 = 0.639
